$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.431.91"
$ws.Range("E2").Value = "  -1.08%  "

$ws.Range("D3").Value = "2.361.78"
$ws.Range("E3").Value = "  +5.25%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.78"
$ws.Range("E5").Value = "  +0.82%  "

$ws.Range("E6").Value = "  +1.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "68.90"
$ws.Range("E7").Value = "  +8.26%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.456"
$ws.Range("E9").Value = "  +0.35%  "

$ws.Range("E10").Value = "  -0.83%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.51"
$ws.Range("E11").Value = "  -1.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.33"
$ws.Range("E12").Value = "  +0.46%  "

$ws.Range("D13").Value = "2.735.52"
$ws.Range("E13").Value = "  +6.11%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.106"
$ws.Range("E14").Value = "  +0.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.66"
$ws.Range("E15").Value = "  +1.19%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.21"
$ws.Range("E16").Value = "  +1.34%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.846"
$ws.Range("E17").Value = "  +1.50%  "

$ws.Range("D18").Value = "2.363.15"
$ws.Range("E18").Value = "  +5.04%  "

$ws.Range("D19").Value = "43.364.17"
$ws.Range("E19").Value = "  -0.92%  "

$ws.Range("D20").Value = "0.0₃0980"
$ws.Range("E20").Value = "  -1.01%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.26"
$ws.Range("E21").Value = "  +3.36%  "

$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.52"
$ws.Range("E22").Value = "  +1.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "249.83"
$ws.Range("E23").Value = "  +0.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.91"
$ws.Range("E24").Value = "  +18.16%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.45"
$ws.Range("E26").Value = "  +1.83%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.21"
$ws.Range("E27").Value = "  -3.52%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.53"
$ws.Range("E28").Value = "  +7.44%  "

$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.88"
$ws.Range("E29").Value = "  +0.25%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.84"
$ws.Range("E30").Value = "  +0.43%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.54"
$ws.Range("E31").Value = "  +9.40%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.128"
$ws.Range("E32").Value = "  -6.85%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.127"
$ws.Range("E33").Value = "  +1.15%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.99"
$ws.Range("E34").Value = "  +4.34%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0690"
$ws.Range("E35").Value = "  -0.37%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.02"
$ws.Range("E36").Value = "  +2.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.47"
$ws.Range("E37").Value = "  +8.87%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.59"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.62"
$ws.Range("E39").Value = "  -0.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0254"
$ws.Range("E40").Value = "  -0.73%  "

$ws.Range("E41").Value = "  +8.62%  "

$ws.Range("E42").Value = "  +0.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "18.27"
$ws.Range("E43").Value = "  +6.92%  "

$ws.Range("E44").Value = "  +9.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "98.73"
$ws.Range("E45").Value = "  +1.34%  "

$ws.Range("E46").Value = "  +2.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.41"
$ws.Range("E47").Value = "  +1.33%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0949"
$ws.Range("E48").Value = "  +0.17%  "

$ws.Range("D49").Value = "1.445.26"
$ws.Range("E49").Value = "  +1.35%  "

$ws.Range("D50").Value = "2.590.40"
$ws.Range("E50").Value = "  +5.64%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.26"
$ws.Range("E51").Value = "  -2.15%  "
